$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D and E columns so numeric-looking strings are preserved exactly
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = '70.339.57'
$ws.Range("E2").Value = '  +3.84%  '

$ws.Range("D3").Value = '3.532.14'
$ws.Range("E3").Value = '  +2.66%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '602.77'
$ws.Range("E5").Value = '  +3.98%  '

$ws.Range("D6").Value = '172.99'
$ws.Range("E6").Value = '  +4.92%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.611'
$ws.Range("E7").Value = '  +1.64%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.520.02'
$ws.Range("E8").Value = '  +2.53%  '

$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").Value = '0.193'
$ws.Range("E10").Value = '  +5.26%  '

$ws.Range("D11").Value = '7.29'
$ws.Range("E11").Value = '  +8.63%  '

$ws.Range("D12").Value = '0.584'
$ws.Range("E12").Value = '  +3.16%  '

$ws.Range("D13").Value = '46.27'
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("D14").Value = '0.0000275'
$ws.Range("E14").Value = '  +2.54%  '

$ws.Range("D15").Value = '4.126.92'
$ws.Range("E15").Value = '  +3.33%  '

$ws.Range("D16").Value = '8.28'
$ws.Range("E16").Value = '  +0.80%  '

$ws.Range("D17").Value = '606.46'
$ws.Range("E17").Value = '  -0.67%  '

$ws.Range("D18").Value = '3.559.77'
$ws.Range("E18").Value = '  +2.94%  '

$ws.Range("D19").Value = '70.481.19'
$ws.Range("E19").Value = '  +3.94%  '

$ws.Range("E20").Value = '  +1.38%  '

$ws.Range("D21").Value = '17.23'
$ws.Range("E21").Value = '  +1.11%  '

$ws.Range("D22").Value = '0.872'
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("D23").Value = '9.26'
$ws.Range("E23").Value = '  -14.88%  '

$ws.Range("D24").Value = '15.64'
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("D25").Value = '96.22'
$ws.Range("E25").Value = '  +1.50%  '

$ws.Range("D26").Value = '3.72'
$ws.Range("E26").Value = '  +0.57%  '

$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").Value = '2.59'
$ws.Range("E28").Value = '  +1.36%  '

$ws.Range("D29").Value = '33.90'
$ws.Range("E29").Value = '  +6.24%  '

$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = '735.36'
$ws.Range("E30").Value = '  +25.01%  '

$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = '9.01'
$ws.Range("E31").Value = '  +1.09%  '

$ws.Range("D32").Value = '3.04'
$ws.Range("E32").Value = '  +0.88%  '

$ws.Range("D33").Value = '8.17'
$ws.Range("E33").Value = '  -1.40%  '

$ws.Range("D34").Value = '7.00'
$ws.Range("E34").Value = '  +4.06%  '

$ws.Range("E35").Value = '  +1.23%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.100'
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '3.56'
$ws.Range("E37").Value = '  +6.39%  '

$ws.Range("D38").Value = '10.72'
$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("D39").Value = '0.0478'
$ws.Range("E39").Value = '  +11.61%  '

$ws.Range("D40").Value = '56.66'
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("E41").Value = '  +0.42%  '

$ws.Range("D42").Value = '0.142'
$ws.Range("E42").Value = '  +5.98%  '

$ws.Range("D43").Value = '3.354.82'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").Value = '0.315'
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("D45").Value = '0.0₃0695'
$ws.Range("E45").Value = '  +2.93%  '

$ws.Range("D46").Value = '32.45'
$ws.Range("E46").Value = '  +0.71%  '

$ws.Range("D47").Value = '2.90'
$ws.Range("E47").Value = '  +7.26%  '

$ws.Range("D48").Value = '2.57'
$ws.Range("E48").Value = '  +4.44%  '

$ws.Range("D49").Value = '0.129'
$ws.Range("E49").Value = '  +1.99%  '

$ws.Range("D50").Value = '133.91'
$ws.Range("E50").Value = '  +1.35%  '

# Remove the temporary text-format styling so cell style indices match the original (unstyled) cells
$ws.Range("D2:E50").ClearFormats()
